$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '29.475.29'
$ws.Range("E2").Value = '  -0.89%  '
$ws.Range("D3").Value = '1.848.30'
$ws.Range("E3").Value = '  -0.15%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.000'
$ws.Range("E4").Value = '  +0.05%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '243.09'
$ws.Range("E5").Value = '  -0.77%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.6499'
$ws.Range("E6").Value = '  +1.88%  '
$ws.Range("E7").Value = '  +0.06%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '47.80'
$ws.Range("E8").Value = '  +3.65%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07496'
$ws.Range("E9").Value = '  -0.11%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.2974'
$ws.Range("E10").Value = '  -0.77%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '24.45'
$ws.Range("E11").Value = '  +1.46%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.07628'
$ws.Range("E12").Value = '  -0.66%  '
$ws.Range("D13").Value = '1.865.63'
$ws.Range("E13").Value = '  +1.51%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.021'
$ws.Range("E14").Value = '  -0.77%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.6845'
$ws.Range("E15").Value = '  -0.45%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '83.52'
$ws.Range("E16").Value = '  -1.15%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.000009449'
$ws.Range("E17").Value = '  -1.51%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '6.101'
$ws.Range("E18").Value = '  +0.17%  '
$ws.Range("D19").Value = '29.539.21'
$ws.Range("E19").Value = '  -0.58%  '
$ws.Range("D20").Value = '2.106.02'
$ws.Range("E20").Value = '  +0.82%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '237.48'
$ws.Range("E21").Value = '  -0.73%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '12.58'
$ws.Range("E22").Value = '  -0.38%  '
$ws.Range("E23").Value = '  +0.06%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '7.705'
$ws.Range("E24").Value = '  +4.94%  '
$ws.Range("E25").Value = '  -0.02%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '157.42'
$ws.Range("E26").Value = '  -1.47%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.1415'
$ws.Range("E27").Value = '  -0.33%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '8.499'
$ws.Range("E28").Value = '  -0.43%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '17.81'
$ws.Range("E29").Value = '  -0.77%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.06075'
$ws.Range("E30").Value = '  +0.65%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.489'
$ws.Range("E31").Value = '  -0.82%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.257'
$ws.Range("E32").Value = '  -0.37%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.138'
$ws.Range("E33").Value = '  -0.02%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.067'
$ws.Range("E34").Value = '  -1.85%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.183'
$ws.Range("E35").Value = '  +2.81%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.856'
$ws.Range("E36").Value = '  -0.69%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.7263'
$ws.Range("E37").Value = '  -1.10%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.594'
$ws.Range("E38").Value = '  -0.68%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.797'
$ws.Range("E39").Value = '  -2.30%  '
$ws.Range("E40").Value = '  +0.08%  '
$ws.Range("D41").Value = '1.200.80'
$ws.Range("E41").Value = '  -2.15%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '6.225'
$ws.Range("E42").Value = '  -2.24%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.9081'
$ws.Range("E43").Value = '  -1.06%  '
$ws.Range("E44").Value = '  -0.06%  '
$ws.Range("D45").Value = '2.015.71'
$ws.Range("E45").Value = '  +0.48%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '101.59'
$ws.Range("E46").Value = '  -0.73%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '66.48'
$ws.Range("E47").Value = '  +0.00%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '7.426'
$ws.Range("E48").Value = '  +10.71%  '
$ws.Range("E49").Value = '  +2.34%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.4046'
$ws.Range("E50").Value = '  -0.93%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '9.150'
$ws.Range("E51").Value = '  -1.94%  '
